# Apply the "456a3b4" data refresh to both the "展览" (sheet 1) and the
# "全部类型" (sheet 4) worksheets, which carry duplicate copies of the
# same table in this workbook.

$wb = $excel.ActiveWorkbook

$sheetIndexes = @(1, 4)

foreach ($sIdx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sIdx)

    # Row 3: "想去人数" (want-to-go count) ticked up by one.
    $ws.Range("F3").Value = 596

    # Row 4 used to be the cancelled 吉安 event; it is dropped and the
    # 鹰潭 event (formerly row 5) takes its place with refreshed numbers.
    $ws.Range("C4").Value = "鹰潭·原神&崩铁&崩坏only"
    $ws.Range("D4").Value = "站江路25号(鹰潭火车站对面) 鹰潭华盛大酒店"
    $ws.Range("F4").Value = 59
    $ws.Range("G4").Value = 65
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=80590"
    $ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202401/HCgQUe0P1704705130296.png"

    # Row 5: now holds the 景德镇 event (formerly row 6).
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("B5").Value = "2024-01-30"
    $ws.Range("C5").Value = "景德镇·原神X崩铁X崩坏动漫展only"
    $ws.Range("D5").Value = "陶阳南路188号 晨枫臻品酒店"
    $ws.Range("E5").Value = "2024.01.30 10:00-01.30 17:00"
    $ws.Range("F5").Value = 35
    $ws.Range("G5").Value = 55
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=80920"
    $ws.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png"

    # Row 6: now holds the 抚州 event (formerly row 7).
    $ws.Range("B6").NumberFormat = "@"
    $ws.Range("B6").Value = "2024-02-01"
    $ws.Range("C6").Value = "抚州·原神&崩铁&崩坏only"
    $ws.Range("D6").Value = "迎宾大道566号 荣耀国际酒店"
    $ws.Range("E6").Value = "2024.02.01 10:00-02.01 17:00"
    $ws.Range("F6").Value = 17
    $ws.Range("G6").Value = 45
    $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=80636"
    $ws.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202401/mVYKfkkX1704781887641.png"

    # Row 7: brand-new 赣州 (cancelled) event appended at the end of the
    # shifted block.
    $ws.Range("B7").NumberFormat = "@"
    $ws.Range("B7").Value = "2024-02-02"
    $ws.Range("C7").Value = "赣州·第三届半夏动漫展（取消）"
    $ws.Range("D7").Value = "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
    $ws.Range("E7").Value = "2024.02.02 10:00-02.04 17:00"
    $ws.Range("F7").Value = 731
    $ws.Range("G7").Value = "不可售"
    $ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=79973"
    $ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202312/eMehCxbh1702972507887.jpeg"

    # Remaining rows only had their "想去人数" counters refreshed.
    $ws.Range("F12").Value = 433
    $ws.Range("F15").Value = 519
    $ws.Range("F17").Value = 271
    $ws.Range("F19").Value = 320
    $ws.Range("F23").Value = 34
    $ws.Range("F24").Value = 18
    $ws.Range("F26").Value = 768
    $ws.Range("F27").Value = 1345
    $ws.Range("F30").Value = 183
    $ws.Range("F31").Value = 52
    $ws.Range("F32").Value = 148
    $ws.Range("F35").Value = 77
    $ws.Range("F39").Value = 1564
    $ws.Range("F44").Value = 3137
    $ws.Range("F46").Value = 155
    $ws.Range("F47").Value = 818
    $ws.Range("F48").Value = 42
}
